# "Changed and fixed John Murungi"
# Updates the DECEMBER 21 rent statement sheet: records December payments
# (PAID column, "G") for several tenants, renames the vacant R-SHOP-5 "LL"
# placeholder to the new tenant OSCAR (with his December rent paid),
# and logs a lump-sum commission payment "PAID ON 15/12".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DECEMBER 21")

# --- Rent "PAID" (column G) entries for December ---------------------------
# R-SHOP-3 / ROSELYDAN LUKWISA
$ws.Range("G7").Value = 3000
# R-SHOP-4 / GLORIA KAGEHA
$ws.Range("G8").Value = 3000

# R-SHOP-5 was vacant ("LL"); now occupied by OSCAR, who paid December rent
$ws.Range("B9").Value = "OSCAR"
$ws.Range("E9").Value = 3000
$ws.Range("G9").Value = 3000

# R-SHOP-6 / LUCY LAMAI
$ws.Range("G10").Value = 3000
# R-SHOP-7 / DORCAS NGIGE
$ws.Range("G11").Value = 3000
# R-SHOP-9 / ERICK MAINA (partial payment)
$ws.Range("G13").Value = 1500

# --- Commission/payments log (row 24) ---------------------------------------
$ws.Range("B24").Value = "PAID ON 15/12"
$ws.Range("E24").Value = 21600
$ws.Range("G24").Value = "PAID ON 15/12"
$ws.Range("I24").Value = 21600

# --- Column B width tweak (widened to fit tenant names) --------------------
$ws.Columns.Item(2).ColumnWidth = 20.7

# --- Sheet view: scroll position / active selection -------------------------
$ws.Activate()
[void]$ws.Range("J24").Select()

# --- Restore the workbook window from a minimized state --------------------
$wb.Windows.Item(1).WindowState = -4143
